$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 36 - this shifts the existing rows 36..140 down to 37..141,
# carrying their data/formatting with them (matches the diff's observed shift pattern).
$ws.Rows.Item(36).Insert()

# Populate the new row 36 with the new daily price record.
$ws.Range("A36").Value = 6
$ws.Range("B36").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C36").Value = "Metropolitana"
$ws.Range("D36").Value = 45014
$ws.Range("E36").Value = 13
$ws.Range("F36").Value = "Fruta"
$ws.Range("G36").Value = 100104
$ws.Range("H36").Value = "Frutos de pepita"
$ws.Range("I36").Value = 100104003
$ws.Range("J36").Value = "Membrillo"
$ws.Range("K36").Value = "Champion"
$ws.Range("L36").Value = "Primera"
$ws.Range("M36").Value = 39
$ws.Range("N36").Value = 230000
$ws.Range("O36").Value = 250000
$ws.Range("P36").Value = 243846
$ws.Range("Q36").Value = "$/bins (450 kilos)"
$ws.Range("R36").Value = "Región de O'Higgins"
$ws.Range("S36").Value = 542
$ws.Range("T36").Value = 450
